# Add "Mid Paper 1" / "Mid Paper 2" columns (F, G) to both sheets, seeded
# with copies of the existing "Paper 1" / "Paper 2" (D / E) marks — done by
# copying the D/E columns and inserting (shift right) so the new cells pick
# up the exact same cell style as the source, instead of a freshly derived
# one. Also makes "Senior Six" the active sheet/tab (selecting E13 there)
# while "Senior Five" keeps F1 selected and is no longer the active tab.

$xlShiftToRight = -4161

$wb = $excel.ActiveWorkbook

$sheetInfo = @(
    @{ Name = "Senior Six";  LastRow = 6; ClearThroughRow = 17 },
    @{ Name = "Senior Five"; LastRow = 6; ClearThroughRow = 6 }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)
    $lastRow = $info.LastRow

    # Copy "Paper 1" (D) into a new "Mid Paper 1" column (F) ...
    $ws.Range("D1:D$lastRow").Copy()
    $ws.Range("F1").Insert($xlShiftToRight)

    # ... and "Paper 2" (E) into a new "Mid Paper 2" column (G).
    $ws.Range("E1:E$lastRow").Copy()
    $ws.Range("G1").Insert($xlShiftToRight)

    # The inserted columns carried the old header labels along with the
    # marks - rename them.
    $ws.Range("F1").Value = "Mid Paper 1"
    $ws.Range("G1").Value = "Mid Paper 2"

    # Insert() shifted/styled the *whole* column (all the way down to
    # whatever the sheet's used range already reached), so undo that
    # below the real data rows.
    if ($info.ClearThroughRow -gt $lastRow) {
        $ws.Range("F" + ($lastRow + 1) + ":G" + $info.ClearThroughRow).Clear()
    }
}

# Senior Five: column C width update + selection moves to F1, tab no
# longer selected.
$wsFive = $wb.Worksheets.Item("Senior Five")
$wsFive.Columns.Item(3).ColumnWidth = 20.46
$wsFive.Range("F1").Select()

# Senior Six becomes the active sheet/tab, with E13 selected.
$wsSix = $wb.Worksheets.Item("Senior Six")
$wsSix.Activate()
$wsSix.Range("E13").Select()
